$wb = $excel.ActiveWorkbook

# Sheet 1: departements
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(7, 3).Value = 44.82758620689656
$ws1.Cells.Item(7, 4).Value = 26
$ws1.Cells.Item(8, 3).Value = 18.75
$ws1.Cells.Item(8, 4).Value = 9
$ws1.Cells.Item(8, 5).Value = 48
$ws1.Cells.Item(16, 3).Value = 5.128205128205128
$ws1.Cells.Item(16, 4).Value = 2
$ws1.Cells.Item(21, 3).Value = 26.08695652173913
$ws1.Cells.Item(21, 4).Value = 18
$ws1.Cells.Item(32, 3).Value = 24.32432432432433
$ws1.Cells.Item(32, 4).Value = 18
$ws1.Cells.Item(44, 3).Value = 21.91780821917808
$ws1.Cells.Item(44, 4).Value = 16
$ws1.Cells.Item(45, 3).Value = 16.66666666666666
$ws1.Cells.Item(45, 4).Value = 5
$ws1.Cells.Item(58, 3).Value = 22.36842105263158
$ws1.Cells.Item(58, 4).Value = 17
$ws1.Cells.Item(61, 3).Value = 45.87155963302752
$ws1.Cells.Item(61, 5).Value = 109
$ws1.Cells.Item(62, 3).Value = 30.76923076923077
$ws1.Cells.Item(62, 4).Value = 20
$ws1.Cells.Item(64, 3).Value = 27.27272727272727
$ws1.Cells.Item(64, 4).Value = 30
$ws1.Cells.Item(64, 5).Value = 110
$ws1.Cells.Item(65, 3).Value = 35.48387096774194
$ws1.Cells.Item(65, 4).Value = 22
$ws1.Cells.Item(76, 3).Value = 29.16666666666667
$ws1.Cells.Item(76, 4).Value = 21
$ws1.Cells.Item(76, 5).Value = 72
$ws1.Cells.Item(80, 3).Value = 51.51515151515152
$ws1.Cells.Item(80, 4).Value = 34
$ws1.Cells.Item(88, 3).Value = 27.27272727272727
$ws1.Cells.Item(88, 5).Value = 33
$ws1.Cells.Item(93, 3).Value = 43.10344827586206
$ws1.Cells.Item(93, 5).Value = 58
$ws1.Cells.Item(96, 3).Value = 62.5
$ws1.Cells.Item(96, 4).Value = 20
$ws1.Cells.Item(99, 3).Value = 11.76470588235294
$ws1.Cells.Item(99, 4).Value = 10
$ws1.Cells.Item(101, 3).Value = 9.433962264150944
$ws1.Cells.Item(101, 4).Value = 5
$ws1.Cells.Item(104, 3).Value = 9.375
$ws1.Cells.Item(104, 4).Value = 12
$ws1.Cells.Item(105, 3).Value = 12.76595744680851
$ws1.Cells.Item(105, 4).Value = 6
$ws1.Cells.Item(106, 3).Value = 5.714285714285714
$ws1.Cells.Item(106, 4).Value = 2
$ws1.Cells.Item(116, 3).Value = 15.55555555555556
$ws1.Cells.Item(116, 4).Value = 7
$ws1.Cells.Item(118, 3).Value = 6.666666666666667
$ws1.Cells.Item(118, 4).Value = 6
$ws1.Cells.Item(129, 3).Value = 7.476635514018691
$ws1.Cells.Item(129, 4).Value = 8
$ws1.Cells.Item(130, 3).Value = 5.88235294117647
$ws1.Cells.Item(130, 4).Value = 10
$ws1.Cells.Item(133, 3).Value = 7.142857142857142
$ws1.Cells.Item(133, 4).Value = 10
$ws1.Cells.Item(137, 3).Value = 6.578947368421052
$ws1.Cells.Item(137, 4).Value = 10
$ws1.Cells.Item(144, 3).Value = 5.376344086021505
$ws1.Cells.Item(144, 4).Value = 5
$ws1.Cells.Item(153, 3).Value = 6.944444444444445
$ws1.Cells.Item(153, 4).Value = 5
$ws1.Cells.Item(155, 3).Value = 4.716981132075472
$ws1.Cells.Item(155, 4).Value = 5
$ws1.Cells.Item(156, 3).Value = 4.545454545454546
$ws1.Cells.Item(156, 4).Value = 5
$ws1.Cells.Item(158, 3).Value = 5.687203791469194
$ws1.Cells.Item(158, 4).Value = 12
$ws1.Cells.Item(159, 3).Value = 6.796116504854369
$ws1.Cells.Item(159, 4).Value = 7
$ws1.Cells.Item(161, 3).Value = 7.428571428571429
$ws1.Cells.Item(161, 4).Value = 13
$ws1.Cells.Item(162, 3).Value = 3.157894736842105
$ws1.Cells.Item(162, 5).Value = 95
$ws1.Cells.Item(163, 3).Value = 9.782608695652174
$ws1.Cells.Item(163, 4).Value = 9
$ws1.Cells.Item(167, 3).Value = 1.075268817204301
$ws1.Cells.Item(167, 5).Value = 93
$ws1.Cells.Item(168, 3).Value = 8.045977011494253
$ws1.Cells.Item(168, 4).Value = 14
$ws1.Cells.Item(169, 3).Value = 11.11111111111111
$ws1.Cells.Item(169, 4).Value = 5
$ws1.Cells.Item(169, 5).Value = 45
$ws1.Cells.Item(170, 3).Value = 9.574468085106384
$ws1.Cells.Item(170, 4).Value = 9
$ws1.Cells.Item(173, 3).Value = 11.23595505617977
$ws1.Cells.Item(173, 4).Value = 10
$ws1.Cells.Item(174, 3).Value = 16.66666666666666
$ws1.Cells.Item(174, 4).Value = 7
$ws1.Cells.Item(177, 3).Value = 5.263157894736842
$ws1.Cells.Item(177, 4).Value = 7
$ws1.Cells.Item(178, 3).Value = 3.92156862745098
$ws1.Cells.Item(178, 4).Value = 2
$ws1.Cells.Item(184, 3).Value = 6.666666666666667
$ws1.Cells.Item(184, 5).Value = 90
$ws1.Cells.Item(187, 3).Value = 8.571428571428571
$ws1.Cells.Item(187, 4).Value = 6
$ws1.Cells.Item(188, 3).Value = 7.272727272727272
$ws1.Cells.Item(188, 4).Value = 4
$ws1.Cells.Item(190, 3).Value = 7.258064516129033
$ws1.Cells.Item(190, 4).Value = 9
$ws1.Cells.Item(191, 3).Value = 3.846153846153846
$ws1.Cells.Item(191, 4).Value = 3
$ws1.Cells.Item(192, 3).Value = 6.976744186046512
$ws1.Cells.Item(192, 4).Value = 6
$ws1.Cells.Item(193, 3).Value = 5
$ws1.Cells.Item(193, 4).Value = 4
$ws1.Cells.Item(194, 3).Value = 3.773584905660377
$ws1.Cells.Item(194, 4).Value = 4
$ws1.Cells.Item(196, 3).Value = 5.263157894736842
$ws1.Cells.Item(196, 4).Value = 5
$ws1.Cells.Item(197, 3).Value = 6.493506493506493
$ws1.Cells.Item(197, 4).Value = 5
$ws1.Cells.Item(198, 3).Value = 7.042253521126761
$ws1.Cells.Item(198, 4).Value = 5
$ws1.Cells.Item(201, 3).Value = 9.48905109489051
$ws1.Cells.Item(201, 4).Value = 13
$ws1.Cells.Item(202, 3).Value = 9.090909090909092
$ws1.Cells.Item(202, 4).Value = 5
$ws1.Cells.Item(202, 5).Value = 55
$ws1.Cells.Item(204, 3).Value = 6.666666666666667
$ws1.Cells.Item(204, 4).Value = 2
$ws1.Cells.Item(208, 3).Value = 7.755102040816326
$ws1.Cells.Item(208, 4).Value = 19
$ws1.Cells.Item(209, 3).Value = 2.912621359223301
$ws1.Cells.Item(209, 4).Value = 3
$ws1.Cells.Item(210, 3).Value = 5.714285714285714
$ws1.Cells.Item(210, 4).Value = 2
$ws1.Cells.Item(222, 3).Value = 2.777777777777778
$ws1.Cells.Item(222, 4).Value = 2
$ws1.Cells.Item(223, 3).Value = 4.098360655737705
$ws1.Cells.Item(223, 4).Value = 5
$ws1.Cells.Item(226, 3).Value = 9.174311926605505
$ws1.Cells.Item(226, 4).Value = 10
$ws1.Cells.Item(238, 3).Value = 8.421052631578947
$ws1.Cells.Item(238, 4).Value = 8
$ws1.Cells.Item(240, 3).Value = 3.08641975308642
$ws1.Cells.Item(240, 4).Value = 5
$ws1.Cells.Item(241, 3).Value = 3.80952380952381
$ws1.Cells.Item(241, 4).Value = 4
$ws1.Cells.Item(252, 3).Value = 1.754385964912281
$ws1.Cells.Item(252, 4).Value = 2
$ws1.Cells.Item(255, 3).Value = 4.273504273504273
$ws1.Cells.Item(255, 4).Value = 10
$ws1.Cells.Item(256, 3).Value = 4.716981132075472
$ws1.Cells.Item(256, 4).Value = 5
$ws1.Cells.Item(258, 3).Value = 5.633802816901409
$ws1.Cells.Item(258, 5).Value = 213
$ws1.Cells.Item(259, 3).Value = 7.547169811320755
$ws1.Cells.Item(259, 4).Value = 8
$ws1.Cells.Item(259, 5).Value = 106
$ws1.Cells.Item(264, 3).Value = 1.96078431372549
$ws1.Cells.Item(264, 5).Value = 102
$ws1.Cells.Item(265, 3).Value = 5.357142857142857
$ws1.Cells.Item(265, 4).Value = 9
$ws1.Cells.Item(266, 3).Value = 11.62790697674419
$ws1.Cells.Item(266, 4).Value = 5
$ws1.Cells.Item(266, 5).Value = 43
$ws1.Cells.Item(267, 3).Value = 5.88235294117647
$ws1.Cells.Item(267, 4).Value = 7
$ws1.Cells.Item(270, 3).Value = 9.433962264150944
$ws1.Cells.Item(270, 4).Value = 10
$ws1.Cells.Item(271, 3).Value = 18
$ws1.Cells.Item(271, 4).Value = 9
$ws1.Cells.Item(272, 3).Value = 7.333333333333333
$ws1.Cells.Item(272, 4).Value = 11
$ws1.Cells.Item(273, 3).Value = 5.389221556886228
$ws1.Cells.Item(273, 4).Value = 9
$ws1.Cells.Item(274, 3).Value = 8.695652173913043
$ws1.Cells.Item(274, 4).Value = 12
$ws1.Cells.Item(276, 3).Value = 1.030927835051546
$ws1.Cells.Item(276, 4).Value = 1
$ws1.Cells.Item(278, 3).Value = 4
$ws1.Cells.Item(278, 4).Value = 2
$ws1.Cells.Item(279, 3).Value = 5.521472392638037
$ws1.Cells.Item(279, 4).Value = 9
$ws1.Cells.Item(281, 3).Value = 5.102040816326531
$ws1.Cells.Item(281, 5).Value = 98
$ws1.Cells.Item(284, 3).Value = 9.210526315789473
$ws1.Cells.Item(284, 4).Value = 7
$ws1.Cells.Item(285, 3).Value = 7.142857142857142
$ws1.Cells.Item(285, 4).Value = 5
$ws1.Cells.Item(287, 3).Value = 9.923664122137405
$ws1.Cells.Item(287, 4).Value = 13
$ws1.Cells.Item(288, 3).Value = 10.12658227848101
$ws1.Cells.Item(288, 4).Value = 8
$ws1.Cells.Item(289, 3).Value = 10.58823529411765
$ws1.Cells.Item(289, 4).Value = 9
$ws1.Cells.Item(290, 3).Value = 9.75609756097561
$ws1.Cells.Item(290, 4).Value = 8
$ws1.Cells.Item(291, 3).Value = 6.422018348623854
$ws1.Cells.Item(291, 4).Value = 7
$ws1.Cells.Item(294, 3).Value = 4.395604395604396
$ws1.Cells.Item(294, 4).Value = 4
$ws1.Cells.Item(305, 3).Value = 0.7246376811594203
$ws1.Cells.Item(305, 4).Value = 2
$ws1.Cells.Item(310, 3).Value = 1.639344262295082
$ws1.Cells.Item(310, 4).Value = 1
$ws1.Cells.Item(312, 3).Value = 0.8130081300813009
$ws1.Cells.Item(312, 4).Value = 1
$ws1.Cells.Item(315, 3).Value = 1.052631578947368
$ws1.Cells.Item(315, 4).Value = 1
$ws1.Cells.Item(316, 3).Value = 2
$ws1.Cells.Item(316, 4).Value = 2
$ws1.Cells.Item(323, 3).Value = 1.481481481481482
$ws1.Cells.Item(323, 4).Value = 2
$ws1.Cells.Item(326, 3).Value = 1.769911504424779
$ws1.Cells.Item(326, 4).Value = 4
$ws1.Cells.Item(328, 3).Value = 1.351351351351351
$ws1.Cells.Item(328, 4).Value = 2
$ws1.Cells.Item(328, 5).Value = 148
$ws1.Cells.Item(329, 3).Value = 2.040816326530612
$ws1.Cells.Item(329, 4).Value = 1
$ws1.Cells.Item(334, 3).Value = 1.587301587301587
$ws1.Cells.Item(334, 4).Value = 1
$ws1.Cells.Item(335, 3).Value = 0.7751937984496124
$ws1.Cells.Item(335, 4).Value = 1
$ws1.Cells.Item(338, 3).Value = 2.727272727272727
$ws1.Cells.Item(338, 4).Value = 3
$ws1.Cells.Item(339, 3).Value = 10
$ws1.Cells.Item(339, 4).Value = 6
$ws1.Cells.Item(343, 3).Value = 1.680672268907563
$ws1.Cells.Item(343, 4).Value = 2
$ws1.Cells.Item(347, 3).Value = 1.063829787234043
$ws1.Cells.Item(347, 4).Value = 1
$ws1.Cells.Item(351, 3).Value = 5
$ws1.Cells.Item(351, 4).Value = 3
$ws1.Cells.Item(352, 3).Value = 2.008032128514056
$ws1.Cells.Item(352, 4).Value = 5
$ws1.Cells.Item(352, 5).Value = 249
$ws1.Cells.Item(355, 3).Value = 1.739130434782609
$ws1.Cells.Item(355, 4).Value = 4
$ws1.Cells.Item(356, 3).Value = 3.333333333333333
$ws1.Cells.Item(356, 4).Value = 4
$ws1.Cells.Item(356, 5).Value = 120
$ws1.Cells.Item(362, 3).Value = 1.515151515151515
$ws1.Cells.Item(362, 4).Value = 3
$ws1.Cells.Item(363, 3).Value = 3.389830508474576
$ws1.Cells.Item(363, 4).Value = 2
$ws1.Cells.Item(363, 5).Value = 59
$ws1.Cells.Item(364, 3).Value = 3.571428571428571
$ws1.Cells.Item(364, 4).Value = 5
$ws1.Cells.Item(365, 3).Value = 1.98019801980198
$ws1.Cells.Item(365, 4).Value = 2
$ws1.Cells.Item(366, 3).Value = 1.111111111111111
$ws1.Cells.Item(366, 4).Value = 1
$ws1.Cells.Item(367, 3).Value = 3.278688524590164
$ws1.Cells.Item(367, 4).Value = 4
$ws1.Cells.Item(369, 3).Value = 2.409638554216868
$ws1.Cells.Item(369, 4).Value = 4
$ws1.Cells.Item(371, 3).Value = 1.36986301369863
$ws1.Cells.Item(371, 4).Value = 2
$ws1.Cells.Item(372, 3).Value = 1.428571428571429
$ws1.Cells.Item(372, 4).Value = 1
$ws1.Cells.Item(376, 3).Value = 1.092896174863388
$ws1.Cells.Item(376, 4).Value = 2
$ws1.Cells.Item(379, 3).Value = 1.351351351351351
$ws1.Cells.Item(379, 5).Value = 74
$ws1.Cells.Item(381, 3).Value = 3.488372093023256
$ws1.Cells.Item(381, 4).Value = 3
$ws1.Cells.Item(382, 3).Value = 1.204819277108434
$ws1.Cells.Item(382, 4).Value = 1
$ws1.Cells.Item(387, 3).Value = 4.705882352941177
$ws1.Cells.Item(387, 4).Value = 4
$ws1.Cells.Item(390, 3).Value = 59.09090909090909
$ws1.Cells.Item(390, 4).Value = 13
$ws1.Cells.Item(395, 3).Value = 51.61290322580645
$ws1.Cells.Item(395, 4).Value = 16
$ws1.Cells.Item(403, 3).Value = 55.88235294117647
$ws1.Cells.Item(403, 4).Value = 19
$ws1.Cells.Item(405, 3).Value = 50
$ws1.Cells.Item(405, 4).Value = 6
$ws1.Cells.Item(417, 3).Value = 35.71428571428572
$ws1.Cells.Item(417, 4).Value = 10
$ws1.Cells.Item(427, 3).Value = 36
$ws1.Cells.Item(427, 4).Value = 9
$ws1.Cells.Item(446, 3).Value = 33.33333333333333
$ws1.Cells.Item(446, 4).Value = 8
$ws1.Cells.Item(452, 3).Value = 47.27272727272727
$ws1.Cells.Item(452, 4).Value = 26
$ws1.Cells.Item(452, 5).Value = 55
$ws1.Cells.Item(453, 3).Value = 59.09090909090909
$ws1.Cells.Item(453, 5).Value = 22
$ws1.Cells.Item(468, 3).Value = 47.91666666666667
$ws1.Cells.Item(468, 4).Value = 23
$ws1.Cells.Item(475, 3).Value = 36.36363636363637
$ws1.Cells.Item(475, 4).Value = 12
$ws1.Cells.Item(475, 5).Value = 33
$ws1.Cells.Item(485, 3).Value = 51.42857142857142
$ws1.Cells.Item(485, 4).Value = 18
$ws1.Cells.Item(518, 3).Value = 17.11711711711712
$ws1.Cells.Item(518, 4).Value = 19
$ws1.Cells.Item(522, 3).Value = 15.49295774647887
$ws1.Cells.Item(522, 5).Value = 71
$ws1.Cells.Item(543, 3).Value = 7.547169811320755
$ws1.Cells.Item(543, 4).Value = 4
$ws1.Cells.Item(544, 3).Value = 31.25
$ws1.Cells.Item(544, 4).Value = 15
$ws1.Cells.Item(546, 3).Value = 14.70588235294118
$ws1.Cells.Item(546, 4).Value = 20
$ws1.Cells.Item(549, 3).Value = 24.52830188679245
$ws1.Cells.Item(549, 4).Value = 26
$ws1.Cells.Item(549, 5).Value = 106
$ws1.Cells.Item(550, 3).Value = 22.22222222222222
$ws1.Cells.Item(550, 4).Value = 10
$ws1.Cells.Item(555, 3).Value = 13.20754716981132
$ws1.Cells.Item(555, 5).Value = 53
$ws1.Cells.Item(559, 3).Value = 18.18181818181818
$ws1.Cells.Item(559, 4).Value = 6
$ws1.Cells.Item(563, 3).Value = 16
$ws1.Cells.Item(563, 4).Value = 12
$ws1.Cells.Item(564, 3).Value = 18.51851851851852
$ws1.Cells.Item(564, 4).Value = 15
$ws1.Cells.Item(565, 3).Value = 19.71830985915493
$ws1.Cells.Item(565, 4).Value = 14
$ws1.Cells.Item(571, 3).Value = 7.692307692307693
$ws1.Cells.Item(571, 4).Value = 5
$ws1.Cells.Item(577, 3).Value = 50
$ws1.Cells.Item(577, 4).Value = 4
$ws1.Cells.Item(577, 5).Value = 8
$ws1.Cells.Item(578, 3).Value = 16.12903225806452
$ws1.Cells.Item(578, 4).Value = 10
$ws1.Cells.Item(578, 5).Value = 62
$ws1.Cells.Item(586, 3).Value = 2.531645569620253
$ws1.Cells.Item(586, 4).Value = 2
$ws1.Cells.Item(589, 3).Value = 3.472222222222222
$ws1.Cells.Item(589, 4).Value = 5
$ws1.Cells.Item(590, 3).Value = 1.408450704225352
$ws1.Cells.Item(590, 4).Value = 1
$ws1.Cells.Item(601, 3).Value = 1.639344262295082
$ws1.Cells.Item(601, 4).Value = 1
$ws1.Cells.Item(603, 3).Value = 1.626016260162602
$ws1.Cells.Item(603, 4).Value = 2
$ws1.Cells.Item(606, 3).Value = 0
$ws1.Cells.Item(606, 4).Value = 0
$ws1.Cells.Item(614, 3).Value = 1.492537313432836
$ws1.Cells.Item(614, 4).Value = 2
$ws1.Cells.Item(615, 3).Value = 2.512562814070352
$ws1.Cells.Item(615, 4).Value = 5
$ws1.Cells.Item(626, 3).Value = 0.7874015748031495
$ws1.Cells.Item(626, 4).Value = 1
$ws1.Cells.Item(629, 3).Value = 0.9259259259259258
$ws1.Cells.Item(629, 4).Value = 1
$ws1.Cells.Item(640, 3).Value = 2.158273381294964
$ws1.Cells.Item(640, 4).Value = 3
$ws1.Cells.Item(643, 3).Value = 2.016129032258065
$ws1.Cells.Item(643, 4).Value = 5
$ws1.Cells.Item(644, 3).Value = 1.680672268907563
$ws1.Cells.Item(644, 4).Value = 2
$ws1.Cells.Item(647, 3).Value = 1.666666666666667
$ws1.Cells.Item(647, 4).Value = 2
$ws1.Cells.Item(647, 5).Value = 120
$ws1.Cells.Item(652, 3).Value = 0.8928571428571428
$ws1.Cells.Item(652, 4).Value = 1
$ws1.Cells.Item(653, 3).Value = 3.03030303030303
$ws1.Cells.Item(653, 4).Value = 6
$ws1.Cells.Item(654, 3).Value = 5.357142857142857
$ws1.Cells.Item(654, 4).Value = 3
$ws1.Cells.Item(654, 5).Value = 56
$ws1.Cells.Item(655, 3).Value = 3.597122302158273
$ws1.Cells.Item(655, 4).Value = 5
$ws1.Cells.Item(658, 3).Value = 4.878048780487805
$ws1.Cells.Item(658, 4).Value = 6
$ws1.Cells.Item(661, 3).Value = 2.312138728323699
$ws1.Cells.Item(661, 4).Value = 4
$ws1.Cells.Item(662, 3).Value = 5.47945205479452
$ws1.Cells.Item(662, 4).Value = 8
$ws1.Cells.Item(669, 5).Value = 117
$ws1.Cells.Item(670, 3).Value = 1.351351351351351
$ws1.Cells.Item(670, 5).Value = 74
$ws1.Cells.Item(673, 3).Value = 2.409638554216868
$ws1.Cells.Item(673, 4).Value = 2
$ws1.Cells.Item(675, 3).Value = 6.666666666666667
$ws1.Cells.Item(675, 4).Value = 9
$ws1.Cells.Item(676, 3).Value = 1.265822784810127
$ws1.Cells.Item(676, 4).Value = 1
$ws1.Cells.Item(678, 3).Value = 5.88235294117647
$ws1.Cells.Item(678, 4).Value = 5
$ws1.Cells.Item(679, 3).Value = 2.678571428571428
$ws1.Cells.Item(679, 4).Value = 3
$ws1.Cells.Item(682, 3).Value = 4.395604395604396
$ws1.Cells.Item(682, 4).Value = 4
$ws1.Cells.Item(683, 3).Value = 2.5
$ws1.Cells.Item(683, 4).Value = 2
$ws1.Cells.Item(686, 3).Value = 4.166666666666666
$ws1.Cells.Item(686, 4).Value = 6
$ws1.Cells.Item(693, 3).Value = 1.805054151624549
$ws1.Cells.Item(693, 4).Value = 5
$ws1.Cells.Item(698, 3).Value = 1.639344262295082
$ws1.Cells.Item(698, 4).Value = 1
$ws1.Cells.Item(704, 3).Value = 2.97029702970297
$ws1.Cells.Item(704, 4).Value = 3
$ws1.Cells.Item(711, 3).Value = 2.222222222222222
$ws1.Cells.Item(711, 4).Value = 3
$ws1.Cells.Item(712, 3).Value = 2.512562814070352
$ws1.Cells.Item(712, 4).Value = 5
$ws1.Cells.Item(714, 3).Value = 2.212389380530973
$ws1.Cells.Item(714, 4).Value = 5
$ws1.Cells.Item(716, 3).Value = 2.013422818791946
$ws1.Cells.Item(716, 4).Value = 3
$ws1.Cells.Item(717, 3).Value = 2.040816326530612
$ws1.Cells.Item(717, 4).Value = 1
$ws1.Cells.Item(722, 3).Value = 1.587301587301587
$ws1.Cells.Item(722, 4).Value = 1
$ws1.Cells.Item(723, 3).Value = 0.7751937984496124
$ws1.Cells.Item(723, 4).Value = 1
$ws1.Cells.Item(726, 3).Value = 3.636363636363636
$ws1.Cells.Item(726, 4).Value = 4
$ws1.Cells.Item(727, 3).Value = 10
$ws1.Cells.Item(727, 4).Value = 6
$ws1.Cells.Item(731, 3).Value = 2.521008403361344
$ws1.Cells.Item(731, 4).Value = 3
$ws1.Cells.Item(735, 3).Value = 1.063829787234043
$ws1.Cells.Item(735, 4).Value = 1
$ws1.Cells.Item(737, 3).Value = 2.158273381294964
$ws1.Cells.Item(737, 4).Value = 3
$ws1.Cells.Item(739, 3).Value = 5
$ws1.Cells.Item(739, 4).Value = 3
$ws1.Cells.Item(740, 3).Value = 2.409638554216868
$ws1.Cells.Item(740, 4).Value = 6
$ws1.Cells.Item(741, 3).Value = 1.680672268907563
$ws1.Cells.Item(741, 4).Value = 2
$ws1.Cells.Item(743, 3).Value = 1.739130434782609
$ws1.Cells.Item(743, 4).Value = 4
$ws1.Cells.Item(744, 3).Value = 3.333333333333333
$ws1.Cells.Item(744, 4).Value = 4
$ws1.Cells.Item(744, 5).Value = 120
$ws1.Cells.Item(749, 3).Value = 0.8928571428571428
$ws1.Cells.Item(749, 4).Value = 1
$ws1.Cells.Item(750, 3).Value = 4.020100502512562
$ws1.Cells.Item(750, 4).Value = 8
$ws1.Cells.Item(751, 3).Value = 5.084745762711865
$ws1.Cells.Item(751, 4).Value = 3
$ws1.Cells.Item(751, 5).Value = 59
$ws1.Cells.Item(752, 3).Value = 4.285714285714286
$ws1.Cells.Item(752, 4).Value = 6
$ws1.Cells.Item(753, 3).Value = 1.98019801980198
$ws1.Cells.Item(753, 4).Value = 2
$ws1.Cells.Item(754, 3).Value = 1.111111111111111
$ws1.Cells.Item(754, 4).Value = 1
$ws1.Cells.Item(755, 3).Value = 4.878048780487805
$ws1.Cells.Item(755, 4).Value = 6
$ws1.Cells.Item(757, 3).Value = 2.409638554216868
$ws1.Cells.Item(757, 4).Value = 4
$ws1.Cells.Item(758, 3).Value = 2.298850574712644
$ws1.Cells.Item(758, 4).Value = 4
$ws1.Cells.Item(759, 3).Value = 5.442176870748299
$ws1.Cells.Item(759, 4).Value = 8
$ws1.Cells.Item(760, 3).Value = 1.428571428571429
$ws1.Cells.Item(760, 4).Value = 1
$ws1.Cells.Item(764, 3).Value = 1.639344262295082
$ws1.Cells.Item(764, 4).Value = 3
$ws1.Cells.Item(766, 5).Value = 117
$ws1.Cells.Item(767, 3).Value = 1.351351351351351
$ws1.Cells.Item(767, 5).Value = 74
$ws1.Cells.Item(769, 3).Value = 3.488372093023256
$ws1.Cells.Item(769, 4).Value = 3
$ws1.Cells.Item(770, 3).Value = 2.409638554216868
$ws1.Cells.Item(770, 4).Value = 2
$ws1.Cells.Item(772, 3).Value = 8.148148148148149
$ws1.Cells.Item(772, 4).Value = 11
$ws1.Cells.Item(773, 3).Value = 2.531645569620253
$ws1.Cells.Item(773, 4).Value = 2
$ws1.Cells.Item(775, 3).Value = 7.058823529411764
$ws1.Cells.Item(775, 4).Value = 6
$ws1.Cells.Item(776, 3).Value = 2.678571428571428
$ws1.Cells.Item(776, 4).Value = 3
$ws1.Cells.Item(784, 3).Value = 1.408450704225352
$ws1.Cells.Item(784, 4).Value = 1
$ws1.Cells.Item(795, 3).Value = 1.639344262295082
$ws1.Cells.Item(795, 4).Value = 1
$ws1.Cells.Item(820, 3).Value = 0.7751937984496124
$ws1.Cells.Item(820, 4).Value = 1
$ws1.Cells.Item(841, 3).Value = 1.666666666666667
$ws1.Cells.Item(841, 4).Value = 2
$ws1.Cells.Item(841, 5).Value = 120
$ws1.Cells.Item(848, 3).Value = 3.389830508474576
$ws1.Cells.Item(848, 4).Value = 2
$ws1.Cells.Item(848, 5).Value = 59
$ws1.Cells.Item(849, 3).Value = 2.857142857142857
$ws1.Cells.Item(849, 4).Value = 4
$ws1.Cells.Item(852, 3).Value = 3.252032520325204
$ws1.Cells.Item(852, 4).Value = 4
$ws1.Cells.Item(856, 3).Value = 1.360544217687075
$ws1.Cells.Item(856, 4).Value = 2
$ws1.Cells.Item(863, 5).Value = 117
$ws1.Cells.Item(864, 3).Value = 1.351351351351351
$ws1.Cells.Item(864, 5).Value = 74
$ws1.Cells.Item(867, 3).Value = 1.204819277108434
$ws1.Cells.Item(867, 4).Value = 1
$ws1.Cells.Item(872, 3).Value = 3.529411764705882
$ws1.Cells.Item(872, 4).Value = 3

# Sheet 2: regions
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 4).Value = 21.68
$ws2.Cells.Item(3, 5).Value = 85
$ws2.Cells.Item(3, 6).Value = 392
$ws2.Cells.Item(6, 4).Value = 51.69
$ws2.Cells.Item(6, 5).Value = 183
$ws2.Cells.Item(6, 6).Value = 354
$ws2.Cells.Item(7, 4).Value = 8.92
$ws2.Cells.Item(7, 5).Value = 75
$ws2.Cells.Item(8, 4).Value = 5.17
$ws2.Cells.Item(8, 5).Value = 45
$ws2.Cells.Item(10, 4).Value = 4.49
$ws2.Cells.Item(10, 5).Value = 39
$ws2.Cells.Item(11, 4).Value = 8.640000000000001
$ws2.Cells.Item(11, 5).Value = 31
$ws2.Cells.Item(13, 4).Value = 42.86
$ws2.Cells.Item(13, 5).Value = 51
$ws2.Cells.Item(16, 4).Value = 3.6
$ws2.Cells.Item(16, 5).Value = 15
$ws2.Cells.Item(17, 4).Value = 1.79
$ws2.Cells.Item(17, 5).Value = 8
$ws2.Cells.Item(18, 4).Value = 0.67
$ws2.Cells.Item(18, 5).Value = 3
$ws2.Cells.Item(19, 4).Value = 1.12
$ws2.Cells.Item(19, 5).Value = 5
$ws2.Cells.Item(20, 4).Value = 8.74
$ws2.Cells.Item(20, 5).Value = 41
$ws2.Cells.Item(20, 6).Value = 469
$ws2.Cells.Item(21, 4).Value = 23.26
$ws2.Cells.Item(21, 5).Value = 60
$ws2.Cells.Item(21, 6).Value = 258
$ws2.Cells.Item(23, 4).Value = 2.42
$ws2.Cells.Item(23, 5).Value = 16
$ws2.Cells.Item(23, 6).Value = 660
$ws2.Cells.Item(24, 4).Value = 19.55
$ws2.Cells.Item(24, 5).Value = 78
$ws2.Cells.Item(25, 4).Value = 7.78
$ws2.Cells.Item(25, 5).Value = 42
$ws2.Cells.Item(25, 6).Value = 540
$ws2.Cells.Item(26, 4).Value = 3.33
$ws2.Cells.Item(26, 5).Value = 22
$ws2.Cells.Item(26, 6).Value = 661
$ws2.Cells.Item(27, 4).Value = 1.66
$ws2.Cells.Item(27, 5).Value = 11
$ws2.Cells.Item(27, 6).Value = 661
$ws2.Cells.Item(28, 4).Value = 2.6
$ws2.Cells.Item(28, 5).Value = 17
$ws2.Cells.Item(28, 6).Value = 655
$ws2.Cells.Item(30, 4).Value = 15.42
$ws2.Cells.Item(30, 5).Value = 35
$ws2.Cells.Item(31, 4).Value = 45.45
$ws2.Cells.Item(31, 5).Value = 55
$ws2.Cells.Item(32, 4).Value = 1.25
$ws2.Cells.Item(32, 5).Value = 7
$ws2.Cells.Item(35, 4).Value = 1.6
$ws2.Cells.Item(35, 5).Value = 9
$ws2.Cells.Item(38, 4).Value = 6.97
$ws2.Cells.Item(38, 5).Value = 45
$ws2.Cells.Item(39, 4).Value = 18.86
$ws2.Cells.Item(39, 5).Value = 76
$ws2.Cells.Item(39, 6).Value = 403
$ws2.Cells.Item(40, 4).Value = 49.27
$ws2.Cells.Item(40, 5).Value = 101
$ws2.Cells.Item(40, 6).Value = 205
$ws2.Cells.Item(41, 6).Value = 799
$ws2.Cells.Item(42, 4).Value = 36.05
$ws2.Cells.Item(42, 5).Value = 137
$ws2.Cells.Item(43, 4).Value = 4.54
$ws2.Cells.Item(43, 5).Value = 33
$ws2.Cells.Item(43, 6).Value = 727
$ws2.Cells.Item(44, 4).Value = 2.12
$ws2.Cells.Item(44, 5).Value = 17
$ws2.Cells.Item(46, 4).Value = 1.76
$ws2.Cells.Item(46, 5).Value = 14
$ws2.Cells.Item(47, 4).Value = 4
$ws2.Cells.Item(47, 5).Value = 28
$ws2.Cells.Item(47, 6).Value = 700
$ws2.Cells.Item(48, 4).Value = 23.27
$ws2.Cells.Item(48, 5).Value = 84
$ws2.Cells.Item(48, 6).Value = 361
$ws2.Cells.Item(50, 6).Value = 860
$ws2.Cells.Item(52, 4).Value = 4.18
$ws2.Cells.Item(52, 5).Value = 32
$ws2.Cells.Item(52, 6).Value = 766
$ws2.Cells.Item(53, 4).Value = 1.74
$ws2.Cells.Item(53, 5).Value = 15
$ws2.Cells.Item(55, 4).Value = 1.52
$ws2.Cells.Item(55, 5).Value = 13
$ws2.Cells.Item(56, 4).Value = 6.14
$ws2.Cells.Item(56, 6).Value = 440
$ws2.Cells.Item(57, 4).Value = 18.67
$ws2.Cells.Item(57, 5).Value = 42
$ws2.Cells.Item(58, 4).Value = 26.8
$ws2.Cells.Item(58, 5).Value = 41
$ws2.Cells.Item(58, 6).Value = 153
$ws2.Cells.Item(59, 4).Value = 1.26
$ws2.Cells.Item(59, 5).Value = 7
$ws2.Cells.Item(59, 6).Value = 556
$ws2.Cells.Item(61, 4).Value = 4.07
$ws2.Cells.Item(61, 5).Value = 20
$ws2.Cells.Item(61, 6).Value = 491
$ws2.Cells.Item(62, 4).Value = 1.44
$ws2.Cells.Item(62, 5).Value = 8
$ws2.Cells.Item(62, 6).Value = 556
$ws2.Cells.Item(63, 6).Value = 556
$ws2.Cells.Item(64, 6).Value = 554
$ws2.Cells.Item(65, 4).Value = 4.51
$ws2.Cells.Item(65, 5).Value = 19
$ws2.Cells.Item(66, 4).Value = 13.08
$ws2.Cells.Item(66, 5).Value = 28
$ws2.Cells.Item(66, 6).Value = 214
$ws2.Cells.Item(68, 4).Value = 1.27
$ws2.Cells.Item(68, 5).Value = 7
$ws2.Cells.Item(68, 6).Value = 551
$ws2.Cells.Item(69, 4).Value = 20.32
$ws2.Cells.Item(69, 5).Value = 63
$ws2.Cells.Item(70, 4).Value = 3.66
$ws2.Cells.Item(70, 5).Value = 17
$ws2.Cells.Item(71, 4).Value = 1.81
$ws2.Cells.Item(71, 5).Value = 10
$ws2.Cells.Item(73, 4).Value = 1.27
$ws2.Cells.Item(73, 5).Value = 7
$ws2.Cells.Item(74, 4).Value = 5.75
$ws2.Cells.Item(74, 5).Value = 50
$ws2.Cells.Item(76, 4).Value = 36.08
$ws2.Cells.Item(76, 5).Value = 92
$ws2.Cells.Item(77, 4).Value = 0.89
$ws2.Cells.Item(77, 5).Value = 10
$ws2.Cells.Item(77, 6).Value = 1124
$ws2.Cells.Item(78, 4).Value = 19.6
$ws2.Cells.Item(78, 6).Value = 602
$ws2.Cells.Item(80, 4).Value = 0.98
$ws2.Cells.Item(80, 5).Value = 11
$ws2.Cells.Item(80, 6).Value = 1124
$ws2.Cells.Item(81, 6).Value = 1124
$ws2.Cells.Item(82, 4).Value = 0.45
$ws2.Cells.Item(82, 5).Value = 5
$ws2.Cells.Item(82, 6).Value = 1123
$ws2.Cells.Item(83, 4).Value = 5.67
$ws2.Cells.Item(83, 5).Value = 49
$ws2.Cells.Item(84, 4).Value = 13.37
$ws2.Cells.Item(84, 5).Value = 73
$ws2.Cells.Item(86, 4).Value = 1.18
$ws2.Cells.Item(86, 5).Value = 13
$ws2.Cells.Item(87, 4).Value = 26.06
$ws2.Cells.Item(87, 5).Value = 147
$ws2.Cells.Item(88, 4).Value = 4.58
$ws2.Cells.Item(88, 5).Value = 41
$ws2.Cells.Item(89, 4).Value = 1.63
$ws2.Cells.Item(89, 5).Value = 18
$ws2.Cells.Item(91, 4).Value = 1.28
$ws2.Cells.Item(91, 5).Value = 14
$ws2.Cells.Item(92, 4).Value = 8.199999999999999
$ws2.Cells.Item(92, 5).Value = 83
$ws2.Cells.Item(92, 6).Value = 1012
$ws2.Cells.Item(93, 4).Value = 21.98
$ws2.Cells.Item(93, 5).Value = 120
$ws2.Cells.Item(94, 4).Value = 46.69
$ws2.Cells.Item(94, 5).Value = 141
$ws2.Cells.Item(94, 6).Value = 302
$ws2.Cells.Item(95, 4).Value = 1.35
$ws2.Cells.Item(95, 5).Value = 18
$ws2.Cells.Item(95, 6).Value = 1331
$ws2.Cells.Item(96, 4).Value = 26.06
$ws2.Cells.Item(96, 5).Value = 190
$ws2.Cells.Item(96, 6).Value = 729
$ws2.Cells.Item(97, 4).Value = 6.01
$ws2.Cells.Item(97, 5).Value = 66
$ws2.Cells.Item(97, 6).Value = 1099
$ws2.Cells.Item(98, 4).Value = 2.17
$ws2.Cells.Item(98, 5).Value = 29
$ws2.Cells.Item(98, 6).Value = 1335
$ws2.Cells.Item(99, 4).Value = 0.9
$ws2.Cells.Item(99, 5).Value = 12
$ws2.Cells.Item(99, 6).Value = 1335
$ws2.Cells.Item(100, 4).Value = 1.73
$ws2.Cells.Item(100, 5).Value = 23
$ws2.Cells.Item(100, 6).Value = 1328
$ws2.Cells.Item(101, 4).Value = 4.96
$ws2.Cells.Item(101, 5).Value = 35
$ws2.Cells.Item(102, 4).Value = 14.68
$ws2.Cells.Item(102, 5).Value = 59
$ws2.Cells.Item(103, 4).Value = 40.11
$ws2.Cells.Item(103, 5).Value = 75
$ws2.Cells.Item(105, 4).Value = 36.12
$ws2.Cells.Item(105, 5).Value = 121
$ws2.Cells.Item(106, 4).Value = 7.41
$ws2.Cells.Item(106, 5).Value = 53
$ws2.Cells.Item(107, 4).Value = 2.22
$ws2.Cells.Item(107, 5).Value = 18
$ws2.Cells.Item(109, 4).Value = 1.86
$ws2.Cells.Item(109, 5).Value = 15

# Sheet 3: national
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 2).Value = 6.37
$ws3.Cells.Item(2, 3).Value = 491
$ws3.Cells.Item(2, 4).Value = 7710
$ws3.Cells.Item(3, 2).Value = 18.66
$ws3.Cells.Item(3, 3).Value = 790
$ws3.Cells.Item(3, 4).Value = 4233
$ws3.Cells.Item(4, 2).Value = 42.78
$ws3.Cells.Item(4, 3).Value = 1019
$ws3.Cells.Item(4, 4).Value = 2382
$ws3.Cells.Item(5, 2).Value = 1.51
$ws3.Cells.Item(5, 3).Value = 148
$ws3.Cells.Item(5, 4).Value = 9795
$ws3.Cells.Item(6, 2).Value = 27.69
$ws3.Cells.Item(6, 3).Value = 1401
$ws3.Cells.Item(6, 4).Value = 5059
$ws3.Cells.Item(7, 2).Value = 5.33
$ws3.Cells.Item(7, 3).Value = 448
$ws3.Cells.Item(7, 4).Value = 8402
$ws3.Cells.Item(8, 2).Value = 2.14
$ws3.Cells.Item(8, 3).Value = 210
$ws3.Cells.Item(8, 4).Value = 9810
$ws3.Cells.Item(9, 2).Value = 1.02
$ws3.Cells.Item(9, 3).Value = 100
$ws3.Cells.Item(9, 4).Value = 9810
$ws3.Cells.Item(10, 2).Value = 1.66
$ws3.Cells.Item(10, 3).Value = 162
$ws3.Cells.Item(10, 4).Value = 9771

Write-Host "Applied all changes"
